# New crime data collected — weekly CompStat report refresh (7th Precinct).
# Updates: report header (issue number / week-ending dates) and the full
# crime-complaints stat grid (rows 15-28, 31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text shared strings collapsed to plain text; same
# font/size/color run throughout so visual result is identical).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# ---------------------------------------------------------------------
# Cells that flip between a numeric value and the text placeholders
# ("0" / "***.*") used elsewhere in the sheet for an empty bucket.
# These need their number format forced to Text before the literal
# string is assigned (otherwise "0" would be reinterpreted as 0), and
# then the surrounding look (General format / right-aligned font) is
# restored by copying formats from an already-correct neighbor cell.
# ---------------------------------------------------------------------

# D15, E15: numeric -> text ("0", "***.*")
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)

# C20: numeric -> text ("0")
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("F15").Copy()
$ws.Range("C20").PasteSpecial(-4122)

# D27, E27: numeric -> text ("0", "***.*")
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("F15").Copy()
$ws.Range("D27:E27").PasteSpecial(-4122)

# C22: text ("0") -> numeric
$ws.Range("C22").Value = 1
$ws.Range("D18").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Plain value refreshes (all keep their existing number format/style).
# ---------------------------------------------------------------------
    # Row 15
    $ws.Range("L15").Value = 0
    # Row 16
    $ws.Range("C16").Value = 2
    $ws.Range("D16").Value = 2
    $ws.Range("G16").Value = 7
    $ws.Range("H16").Value = 0
    $ws.Range("I16").Value = 110
    $ws.Range("J16").Value = 95
    $ws.Range("K16").Value = 15.789473684210
    $ws.Range("L16").Value = -12.698412698412
    $ws.Range("M16").Value = -12.698412698412
    $ws.Range("N16").Value = -85.488126649076
    # Row 17
    $ws.Range("C17").Value = 3
    $ws.Range("D17").Value = 5
    $ws.Range("E17").Value = -40
    $ws.Range("F17").Value = 19
    $ws.Range("H17").Value = 18.75
    $ws.Range("I17").Value = 198
    $ws.Range("J17").Value = 218
    $ws.Range("K17").Value = -9.174311926605
    $ws.Range("L17").Value = 1.538461538461
    $ws.Range("M17").Value = 63.636363636363
    $ws.Range("N17").Value = -6.603773584905
    # Row 18
    $ws.Range("C18").Value = 3
    $ws.Range("E18").Value = 200
    $ws.Range("F18").Value = 5
    $ws.Range("G18").Value = 8
    $ws.Range("H18").Value = -37.5
    $ws.Range("I18").Value = 87
    $ws.Range("J18").Value = 103
    $ws.Range("K18").Value = -15.533980582524
    $ws.Range("L18").Value = -45.625
    $ws.Range("M18").Value = -6.451612903225
    $ws.Range("N18").Value = -75.766016713091
    # Row 19
    $ws.Range("C19").Value = 11
    $ws.Range("D19").Value = 2
    $ws.Range("E19").Value = 450
    $ws.Range("G19").Value = 38
    $ws.Range("H19").Value = 13.157894736842
    $ws.Range("I19").Value = 412
    $ws.Range("J19").Value = 474
    $ws.Range("K19").Value = -13.080168776371
    $ws.Range("L19").Value = -24.403669724770
    $ws.Range("M19").Value = 56.653992395437
    $ws.Range("N19").Value = -1.670644391408
    # Row 20
    $ws.Range("E20").Value = -100
    $ws.Range("G20").Value = 3
    $ws.Range("H20").Value = 33.333333333333
    $ws.Range("J20").Value = 51
    $ws.Range("K20").Value = -37.254901960784
    $ws.Range("L20").Value = -47.540983606557
    $ws.Range("M20").Value = -31.914893617021
    $ws.Range("N20").Value = -91.208791208791
    # Row 21
    $ws.Range("C21").Value = 19
    $ws.Range("D21").Value = 11
    $ws.Range("E21").Value = 72.727272727272
    $ws.Range("F21").Value = 78
    $ws.Range("G21").Value = 75
    $ws.Range("H21").Value = 4
    $ws.Range("I21").Value = 851
    $ws.Range("J21").Value = 953
    $ws.Range("K21").Value = -10.703043022035
    $ws.Range("L21").Value = -22.565969062784
    $ws.Range("M21").Value = 28.939393939393
    $ws.Range("N21").Value = -60.215053763440
    # Row 22
    $ws.Range("D22").Value = 2
    $ws.Range("E22").Value = -50
    $ws.Range("G22").Value = 3
    $ws.Range("H22").Value = -33.333333333333
    $ws.Range("I22").Value = 25
    $ws.Range("J22").Value = 21
    $ws.Range("K22").Value = 19.047619047619
    $ws.Range("L22").Value = 38.888888888888
    $ws.Range("M22").Value = 56.25
    # Row 23
    $ws.Range("C23").Value = 2
    $ws.Range("E23").Value = -33.333333333333
    $ws.Range("G23").Value = 18
    $ws.Range("H23").Value = -11.111111111111
    $ws.Range("I23").Value = 160
    $ws.Range("J23").Value = 203
    $ws.Range("K23").Value = -21.182266009852
    $ws.Range("L23").Value = 2.564102564102
    $ws.Range("M23").Value = 29.032258064516
    # Row 24
    $ws.Range("C24").Value = 20
    $ws.Range("D24").Value = 15
    $ws.Range("E24").Value = 33.333333333333
    $ws.Range("F24").Value = 118
    $ws.Range("G24").Value = 114
    $ws.Range("H24").Value = 3.508771929824
    $ws.Range("I24").Value = 1215
    $ws.Range("J24").Value = 1303
    $ws.Range("K24").Value = -6.753645433614
    $ws.Range("L24").Value = 6.113537117903
    $ws.Range("M24").Value = 75.578034682080
    # Row 25
    $ws.Range("C25").Value = 11
    $ws.Range("D25").Value = 10
    $ws.Range("E25").Value = 10
    $ws.Range("F25").Value = 69
    $ws.Range("G25").Value = 67
    $ws.Range("H25").Value = 2.985074626865
    $ws.Range("I25").Value = 709
    $ws.Range("J25").Value = 859
    $ws.Range("K25").Value = -17.462165308498
    $ws.Range("L25").Value = 8.409785932721
    # Row 26
    $ws.Range("C26").Value = 7
    $ws.Range("D26").Value = 8
    $ws.Range("E26").Value = -12.5
    $ws.Range("G26").Value = 37
    $ws.Range("H26").Value = -21.621621621621
    $ws.Range("I26").Value = 381
    $ws.Range("J26").Value = 417
    $ws.Range("K26").Value = -8.633093525179
    $ws.Range("L26").Value = -7.971014492753
    $ws.Range("M26").Value = 25.742574257425
    # Row 27
    $ws.Range("L27").Value = -36.842105263157
    # Row 28
    $ws.Range("C28").Value = 2
    $ws.Range("F28").Value = 5
    $ws.Range("H28").Value = 400
    $ws.Range("I28").Value = 45
    $ws.Range("K28").Value = 28.571428571428
    $ws.Range("L28").Value = 2.272727272727
    # Row 31
    $ws.Range("L31").Value = 100
